$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell G1: "Cargos" ---
$ws.Range("G1").Value = "Cargos"

# --- Column E (new): hours worked per step ---
$ws.Range("E2").Value = 0.020833333333333332
$ws.Range("E2").NumberFormat = "h:mm"

$ws.Range("E3").Value = 0.03125
$ws.Range("E3").NumberFormat = "h:mm"

$ws.Range("E4").Value = 0.010416666666666666
$ws.Range("E4").NumberFormat = "h:mm"

$ws.Range("E5").Value = 0.006944444444444444
$ws.Range("E5").NumberFormat = "h:mm"

$ws.Range("E6").Value = 0.25
$ws.Range("E6").NumberFormat = "h:mm"

# --- Column G (new): "Cargos" values ---
# Create the plain "centered, general format" style first (G3/G4/G5) so it
# lands at cellXfs index 3, then the "centered, h:mm format" style (G2) so
# it lands at cellXfs index 4 - matching the order the styles were authored.
$ws.Range("G3").Value = ""
$ws.Range("G3").HorizontalAlignment = -4108

$ws.Range("G4").Value = ""
$ws.Range("G4").HorizontalAlignment = -4108

$ws.Range("G5").Value = ""
$ws.Range("G5").HorizontalAlignment = -4108

$ws.Range("G2").Value = 0.42708333333333331
$ws.Range("G2").NumberFormat = "h:mm"
$ws.Range("G2").HorizontalAlignment = -4108

$ws.Range("G6").Value = 0.09375
$ws.Range("G6").NumberFormat = "h:mm"

# --- Merge G2:G5 ---
$ws.Range("G2:G5").Merge()

# --- Row 7 totals ---
$ws.Range("E7").Formula = "=SUM(E2:E6)"
$ws.Range("E7").NumberFormat = "h:mm"

$ws.Range("F7").Formula = "=SUM(F2:F6)"
$ws.Range("F7").NumberFormat = "h:mm"

$ws.Range("G7").Formula = "=SUM(G2:G6)"
$ws.Range("G7").NumberFormat = "h:mm"

# --- Column widths ---
$ws.Columns("C").ColumnWidth = 36.5
$ws.Columns("F").ColumnWidth = 4.666666666666667
$ws.Columns("G").ColumnWidth = 6

# --- Selection ---
$ws.Range("E7:G7").Select()
